$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the "Databas" label from A15 up to A14 (no row shifting involved)
$ws.Cells.Item(15, 1).ClearContents()
$ws.Cells.Item(14, 1).Value = "Databas"

# Add merged note cell B14:C14 with centered text
$ws.Range("B14:C14").Merge()
$ws.Cells.Item(14, 2).Value = "(Häma,Skicka,Koppla)"
$ws.Range("B14:C14").HorizontalAlignment = -4108

# Row 16: C16 changes from 15 to 120
$ws.Cells.Item(16, 3).Value = 120

# Row 17: add a new date value in B17 (keep C17 = 15 as-is)
$ws.Cells.Item(16, 2).Copy()
$ws.Cells.Item(17, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(17, 2).Value = 44956

# Update selection to match target
$ws.Range("D23").Select()
